$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from the existing
# Unnamed: 28 header cell so the new headers match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every data row (rows 2-49)
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 61
    $ws.Cells.Item($r, 31).Value = 101
    $ws.Cells.Item($r, 32).Value = 0
}
